$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# "... template coding (HTML, CSS, JS), data integration (SAS), Adobe Campaign
#  workflow setup, and segmentation configuration (Queries)."
#  -> "... data extraction (SAS) ... segmentation setup (Queries)."
Replace-Text ", data integration" ", data extraction"
Replace-Text "workflow setup, and segmentation configuration" "workflow setup, and segmentation setup"

# "Ad-hoc analysis ... improve resources allocation and make recommendations."
# -> "... improve resources allocation, make recommendations and identify actionable insights."
Replace-Text "improve resources allocation and make recommendations." "improve resources allocation, make recommendations and identify actionable insights."

# Skills line: drop "Domo", add "Looker Studio" and "BigQuery"
Replace-Text "Data Analysis and visualization (Domo, Tableau, Pandas, " "Data Analysis and visualization (Tableau, Looker Studio, Pandas, "
Replace-Text "). SQL, Python, SAS." "). SQL, BigQuery, Python, SAS."
